# feat(init): implement ReadClaimsData workflow and business/system exception handling
#
# This script adds the new configuration entries required by the
# ReadClaimsData workflow (SQLite DB path, claims sheet name) as well as
# the email notification settings (sender/recipients, subjects & bodies)
# used by the business- and system-exception handling branches, on the
# "Constants" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Constants" sheet

# --- SQLiteDBPath (new row, inserted right after ClaimsSummaryFilePath) ---
$ws.Range("A21").Value2 = "SQLiteDBPath"
$ws.Range("B21").Value2 = "C:\Users\MorisMwaiWachira\Desktop\MorisMwai_RPA_Assignment\Database\claims.db"

# --- SheetName (appended after MaxClaimAmount) ---
$ws.Range("A27").Value2 = "SheetName"
$ws.Range("B27").Value2 = "ClaimsData"

# --- Email notification settings ---
$ws.Range("A29").Value2 = "SenderEmail"
$ws.Range("B29").Value2 = "mwaimoris@outlook.com"

$ws.Range("A30").Value2 = "RecipientEmail_Business"
$ws.Range("B30").Value2 = "mwaimoris@gmail.com"

$ws.Range("A31").Value2 = "RecipientEmail_System"
$ws.Range("B31").Value2 = "morismwai1@gmail.com"

# --- No file found notification ---
$ws.Range("A33").Value2 = "Subject_NoFileFound"
$ws.Range("B33").Value2 = "Missing Input File"

$ws.Range("A34").Value2 = "Body_NoFileFound"
$bodyNoFileFound = @"
Hello,
The automation was unable to start because the required input file was not found at path:
@InputFilePath
Please upload the file and notify the automation team.
Timesstamp: @Timestamp
Best regards,
Robot :)
"@
$ws.Range("B34").Value2 = $bodyNoFileFound
$ws.Range("B34").WrapText = $true
$ws.Rows.Item(34).RowHeight = 180

# --- Empty claims file notification ---
$ws.Range("A36").Value2 = "Subject_EmptyFile"
$ws.Range("B36").Value2 = "Empty Claims File"

$ws.Range("A37").Value2 = "Body_EmptyFile"
$bodyEmptyFile = @"
Hello, 
The automation could not proceed because the claims data file is empty. 
Please verify the data and reload the file. 
Timestamp: @Timestamp 
Best regards,
Robot :)
"@
$ws.Range("B37").Value2 = $bodyEmptyFile
$ws.Range("B37").WrapText = $true
$ws.Rows.Item(37).RowHeight = 150

# --- Invalid Excel format notification ---
$ws.Range("A39").Value2 = "Subject_InvalidExcelFormat"
$ws.Range("B39").Value2 = "Invalid Excel Format"

$ws.Range("A40").Value2 = "Body_InvalidExcelFormat"
$bodyInvalidExcelFormat = @"
Hello,
The automation was unable to process the file due to missing or incorrect headers in the Excel sheet. 
Please correct the format and upload again. 
Timestamp: @Timestamp 
Best regards, 
Robot :)
"@
$ws.Range("B40").Value2 = $bodyInvalidExcelFormat
$ws.Range("B40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 165

# --- System exception notification ---
$ws.Range("A42").Value2 = "Subject_SystemException"
$ws.Range("B42").Value2 = "Automation Error!"

$ws.Range("A43").Value2 = "Body_SystemException"
$bodySystemException = @"
Hello,
An exception occurred during the automation process. Please find the details below:
Exception Source: @Source
Exception Message: @Message
A screenshot of the error has been attached for reference. Please see the attachment for more details.
Thank you and have a good day,
Robot :)
"@
$ws.Range("B43").Value2 = $bodySystemException
$ws.Range("B43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 180

# --- Widen column B so the new e-mail / subject values are fully visible ---
$ws.Columns.Item(2).AutoFit() | Out-Null

Write-Host "Constants sheet updated with ReadClaimsData + exception-handling configuration."
